$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date value (serial 45203 = 2023-10-04) for every
# data row (rows 2-264). This was bulk-updated to serial 45205 (2023-10-06).
for ($r = 2; $r -le 264; $r++) {
    $ws.Cells.Item($r, 3).Value = 45205
}
